# Applies the "vozrazhenie_na otzyv_otvetchika" template edit:
# replace literal values with {{placeholder}} tokens, update the
# defendant address paragraph (text + font-size bump on the paragraph
# mark), and split three "lawyer_*" runs into {{ / name / }} pieces so
# that the inner token keeps a white highlight.

$d = $word.ActiveDocument

function Replace-AllText($findText, $replaceText) {
    $rng = $d.Content
    while ($rng.Find.Execute($findText, $true, $false, $false, $false, $false,
                              $true, 1, $false, $replaceText, 2)) {
        $rng.Collapse(0)
    }
}

# 1. case number
Replace-AllText "03/0463/0059" "{{legalcase_num}}"

# 2. court name
Replace-AllText "В Арбитражный суд города Москвы" "В {{court_name}}"

# 3. court index/address
Replace-AllText "115225, г. Москва, ул. Большая Тульская, 17 " "{{court_index}}, {{court_address}} "

# 8/12. defendant name (two occurrences)
Replace-AllText 'ООО "РестоБар"' "{{defendant_name}}"

# 9. defendant inn/ogrn
Replace-AllText "ИНН:7703417249 ,ОГРН: 1167746906997" "ИНН:{{defendant_inn}}, ОГРН: {{defendant_ogrn}}"

# 10. defendant address
Replace-AllText "Адрес: 123112, МОСКВА Г, ПРЕСНЕНСКАЯ НАБ, ДОМ 12, ЭТАЖ 75" "Адрес: {{defendant_address}}"

# 11. case number (the A40- case)
Replace-AllText "Дело № А40-289710/2021" "Дело № {{case_num}}"

# 13. violation domain (inside the longer sentence run)
Replace-AllText "ресторан-москва-сити.рф" "{{violation_domain}}"

Write-Host "Simple replacements done"
